# CORS errors and critical HBD transfer bug corrected
#
# - Fill in the missing "ok" results for the Create+top-up / Pay with
#   top-up columns (C & D) across all device rows.
# - Flag the iPhone Safari / Pay with top-up cell as broken ("NOK", red font).
# - Narrow column C slightly so it no longer shares column D's width.
# - Restore the page setup (paper size / orientation) for printing.
# - Selection cursor left on E12, matching the author's last click.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "ok" status cells that were missing -----------------
# Row 2 (desktop): C2 already "ok" -> add D2
$ws.Range("D2").Value = "ok"

# Row 3 (android samsung browser): add C3, D3, F3
$ws.Range("C3").Value = "ok"
$ws.Range("D3").Value = "ok"
$ws.Range("F3").Value = "ok"

# Row 4 (android chrome): add C4, D4
$ws.Range("C4").Value = "ok"
$ws.Range("D4").Value = "ok"

# Row 5 (iPhone Safari): add C5, D5, and the broken E5 cell
$ws.Range("C5").Value = "ok"
$ws.Range("D5").Value = "ok"

# Row 6 (iPhone Chrome): add C6, D6
$ws.Range("C6").Value = "ok"
$ws.Range("D6").Value = "ok"

# --- The critical bug: Pay with top-up is broken on iPhone Safari ----
$ws.Range("E5").Value = "NOK"
$ws.Range("E5").Font.Color = 255

# --- Column C is narrower than before; D keeps the old shared width --
$ws.Columns.Item(3).ColumnWidth = 18.666666666666668

# --- Restore printable page setup ------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Leave the selection where the author last clicked ----------------
$ws.Range("E12").Select() | Out-Null
